# Weekly refresh: reshuffle the per-row price/quality data (columns D, I, J,
# K, L, M, P) for rows 2-20 according to the new weekly ordering. Columns
# A, B, C, E, F, G, H, N, O, Q, R are constant across all rows so they are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "before" values for the columns that move, keyed by row.
$cols = @("D", "I", "J", "K", "L", "M", "P")
$snapshot = @{}
for ($r = 2; $r -le 20; $r++) {
    $row = @{}
    foreach ($c in $cols) {
        $row[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $row
}

# New row -> source (old) row mapping.
$mapping = @{
    2  = 20
    3  = 4
    4  = 5
    5  = 9
    6  = 10
    7  = 3
    8  = 2
    9  = 19
    10 = 7
    11 = 11
    12 = 12
    13 = 13
    14 = 17
    15 = 18
    16 = 8
    17 = 6
    18 = 16
    19 = 14
    20 = 15
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $srcData[$c]
    }
}
